# Apply Behemoth_Profits market-data refresh per commit diff.
# Each (sheet, row) group below corresponds to one changed data row
# in the underlying OOXML diff; columns H-N hold the derived price/
# profit figures that were refreshed by the scheduled runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 863.7646999999999
$ws.Range("I28").Value = 562.2222
$ws.Range("J28").Value = 1203
$ws.Range("K28").Value = 562.2222
$ws.Range("L28").Value = 1203
$ws.Range("M28").Value = -77.22220000000004
$ws.Range("N28").Value = -2173
$ws.Range("H62").Value = 9492.143
$ws.Range("I62").Value = 4855.5
$ws.Range("J62").Value = 13707.272
$ws.Range("K62").Value = 4855.5
$ws.Range("L62").Value = 13707.272
$ws.Range("M62").Value = -4231.5
$ws.Range("N62").Value = -14955.272
$ws.Range("H65").Value = 9492.143
$ws.Range("I65").Value = 4855.5
$ws.Range("J65").Value = 13707.272
$ws.Range("K65").Value = 24277.5
$ws.Range("L65").Value = 68536.36
$ws.Range("M65").Value = -21157.5
$ws.Range("N65").Value = -74776.36
$ws.Range("H93").Value = 100000
$ws.Range("J93").Value = 100000
$ws.Range("L93").Value = 100000
$ws.Range("N93").Value = -104992
$ws.Range("H98").Value = 76924670
$ws.Range("I98").Value = 90910710
$ws.Range("J98").Value = 1450
$ws.Range("K98").Value = 90910710
$ws.Range("L98").Value = 1450
$ws.Range("M98").Value = -90909212
$ws.Range("N98").Value = -4446
$ws.Range("H122").Value = 76924670
$ws.Range("I122").Value = 90910710
$ws.Range("J122").Value = 1450
$ws.Range("K122").Value = 272732130
$ws.Range("L122").Value = 4350
$ws.Range("M122").Value = -272729680
$ws.Range("N122").Value = -9250
$ws.Range("H135").Value = 2483
$ws.Range("I135").Value = 5949.5
$ws.Range("J135").Value = 1327.5
$ws.Range("K135").Value = 53545.5
$ws.Range("L135").Value = 11947.5
$ws.Range("M135").Value = -51010.5
$ws.Range("N135").Value = -17017.5
$ws.Range("H137").Value = 198130.02
$ws.Range("I137").Value = 258184.25
$ws.Range("J137").Value = 2953.75
$ws.Range("K137").Value = 774552.75
$ws.Range("L137").Value = 8861.25
$ws.Range("M137").Value = -772002.75
$ws.Range("N137").Value = -13961.25
$ws.Range("H138").Value = 2564.4783
$ws.Range("J138").Value = 2712.5083
$ws.Range("L138").Value = 8137.5249
$ws.Range("N138").Value = -18417.5249

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6102415.5
$ws.Range("I32").Value = 7354656
$ws.Range("K32").Value = 7354656
$ws.Range("M32").Value = -7354369
$ws.Range("H45").Value = 22728864
$ws.Range("I45").Value = 26317218
$ws.Range("K45").Value = 26317218
$ws.Range("M45").Value = -26316841
$ws.Range("H63").Value = 5606.2666
$ws.Range("J63").Value = 7687.25
$ws.Range("L63").Value = 7687.25
$ws.Range("N63").Value = -9059.25
$ws.Range("H66").Value = 5606.2666
$ws.Range("J66").Value = 7687.25
$ws.Range("L66").Value = 38436.25
$ws.Range("N66").Value = -45300.25
$ws.Range("H74").Value = 4905368
$ws.Range("I74").Value = 6579955.5
$ws.Range("J74").Value = 10419.462
$ws.Range("K74").Value = 6579955.5
$ws.Range("L74").Value = 10419.462
$ws.Range("M74").Value = -6579081.5
$ws.Range("N74").Value = -12167.462
$ws.Range("H77").Value = 4905368
$ws.Range("I77").Value = 6579955.5
$ws.Range("J77").Value = 10419.462
$ws.Range("K77").Value = 32899777.5
$ws.Range("L77").Value = 52097.31
$ws.Range("M77").Value = -32895409.5
$ws.Range("N77").Value = -60833.31
$ws.Range("H102").Value = 27631.166
$ws.Range("I102").Value = 27631.166
$ws.Range("K102").Value = 27631.166
$ws.Range("M102").Value = -26009.166
$ws.Range("H106").Value = 52320
$ws.Range("J106").Value = 52320
$ws.Range("L106").Value = 52320
$ws.Range("N106").Value = -54844
$ws.Range("H132").Value = 3144
$ws.Range("I132").Value = 1869.7916
$ws.Range("J132").Value = 13337.667
$ws.Range("K132").Value = 5609.3748
$ws.Range("L132").Value = 40013.001
$ws.Range("M132").Value = -3079.3748
$ws.Range("N132").Value = -45073.001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1551.44
$ws.Range("I94").Value = 567.6667
$ws.Range("J94").Value = 3027.1
$ws.Range("K94").Value = 567.6667
$ws.Range("L94").Value = 3027.1
$ws.Range("M94").Value = -116.6667
$ws.Range("N94").Value = -3929.1
$ws.Range("H105").Value = 1896.8422
$ws.Range("I105").Value = 1297.0588
$ws.Range("J105").Value = 2382.3809
$ws.Range("K105").Value = 1297.0588
$ws.Range("L105").Value = 2382.3809
$ws.Range("M105").Value = 449.9412
$ws.Range("N105").Value = -5876.3809
$ws.Range("H134").Value = 336098.2
$ws.Range("I134").Value = 1278.5714
$ws.Range("J134").Value = 922032.5
$ws.Range("K134").Value = 3835.7142
$ws.Range("L134").Value = 2766097.5
$ws.Range("M134").Value = -1300.7142
$ws.Range("N134").Value = -2771167.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 817778.2
$ws.Range("I31").Value = 16048.454
$ws.Range("K31").Value = 16048.454
$ws.Range("M31").Value = -15753.454
$ws.Range("H34").Value = 817778.2
$ws.Range("I34").Value = 16048.454
$ws.Range("K34").Value = 16048.454
$ws.Range("M34").Value = -15846.454
$ws.Range("H103").Value = 40144
$ws.Range("I103").Value = 9999.5
$ws.Range("K103").Value = 9999.5
$ws.Range("M103").Value = -8827.5
$ws.Range("H105").Value = 1077.4667
$ws.Range("I105").Value = 1010.7857
$ws.Range("K105").Value = 1010.7857
$ws.Range("M105").Value = 736.2143
$ws.Range("H107").Value = 3112.5454
$ws.Range("I107").Value = 2026.5555
$ws.Range("K107").Value = 2026.5555
$ws.Range("M107").Value = -106.5554999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 949.125
$ws.Range("I23").Value = 1525.25
$ws.Range("J23").Value = 373
$ws.Range("K23").Value = 4575.75
$ws.Range("L23").Value = 1119
$ws.Range("M23").Value = -4340.75
$ws.Range("N23").Value = -1589
$ws.Range("H68").Value = 2888.6667
$ws.Range("J68").Value = 3777.3333
$ws.Range("L68").Value = 11331.9999
$ws.Range("N68").Value = -12953.9999
$ws.Range("H71").Value = 2888.6667
$ws.Range("J71").Value = 3777.3333
$ws.Range("L71").Value = 33995.9997
$ws.Range("N71").Value = -42107.9997
$ws.Range("H113").Value = 1337.6316
$ws.Range("I113").Value = 714
$ws.Range("K113").Value = 2142
$ws.Range("M113").Value = 28
$ws.Range("H131").Value = 18032.555
$ws.Range("J131").Value = 35044.332
$ws.Range("L131").Value = 105132.996
$ws.Range("N131").Value = -115212.996

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 16666.445
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 16666.445
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 16666.445
$ws.Range("M43").ClearContents()
$ws.Range("N43").Value = -16968.445
$ws.Range("H46").Value = 30000
$ws.Range("J46").Value = 30000
$ws.Range("L46").Value = 30000
$ws.Range("N46").Value = -30312
$ws.Range("H57").Value = 10000
$ws.Range("J57").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("N57").ClearContents()
$ws.Range("H80").Value = 5554.25
$ws.Range("I80").Value = 4132.4
$ws.Range("J80").Value = 7924
$ws.Range("K80").Value = 4132.4
$ws.Range("L80").Value = 7924
$ws.Range("M80").Value = -3134.4
$ws.Range("N80").Value = -9920
$ws.Range("H83").Value = 5554.25
$ws.Range("I83").Value = 4132.4
$ws.Range("J83").Value = 7924
$ws.Range("K83").Value = 20662
$ws.Range("L83").Value = 39620
$ws.Range("M83").Value = -15670
$ws.Range("N83").Value = -49604
$ws.Range("H102").Value = 2060.4482
$ws.Range("I102").Value = 1269.3334
$ws.Range("J102").Value = 2908.0715
$ws.Range("K102").Value = 1269.3334
$ws.Range("L102").Value = 2908.0715
$ws.Range("M102").Value = 352.6666
$ws.Range("N102").Value = -6152.0715
$ws.Range("H132").Value = 22729598
$ws.Range("J132").Value = 1680.1111
$ws.Range("L132").Value = 5040.3333
$ws.Range("N132").Value = -10100.3333

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3069.0322
$ws.Range("I40").Value = 2657.739
$ws.Range("K40").Value = 2657.739
$ws.Range("M40").Value = -2521.739
$ws.Range("H112").Value = 110000
$ws.Range("J112").Value = 110000
$ws.Range("L112").Value = 110000
$ws.Range("N112").Value = -112954
$ws.Range("H122").Value = 5708.5
$ws.Range("I122").Value = 5174.25
$ws.Range("K122").Value = 15522.75
$ws.Range("M122").Value = -13072.75
$ws.Range("H132").Value = 365105.6
$ws.Range("I132").Value = 9258.825999999999
$ws.Range("K132").Value = 27776.478
$ws.Range("M132").Value = -25246.478

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 151000
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("N4").ClearContents()
